$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 212; this shifts existing rows 212-259 down to 213-260
$ws.Rows.Item(212).Insert()

# Populate the newly inserted row 212 with the new weekly record
$ws.Cells.Item(212, 1).Value = 8
$ws.Cells.Item(212, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(212, 3).Value = "Coquimbo"
$ws.Cells.Item(212, 4).Value = 44782
$ws.Cells.Item(212, 5).Value = 4
$ws.Cells.Item(212, 6).Value = 100112031
$ws.Cells.Item(212, 7).Value = "Poroto verde"
$ws.Cells.Item(212, 8).Value = "Magnum"
$ws.Cells.Item(212, 9).Value = "Primera"
$ws.Cells.Item(212, 10).Value = 520
$ws.Cells.Item(212, 11).Value = 34000
$ws.Cells.Item(212, 12).Value = 35000
$ws.Cells.Item(212, 13).Value = 34500
$ws.Cells.Item(212, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(212, 15).Value = "Perú"
$ws.Cells.Item(212, 16).Value = 1380
$ws.Cells.Item(212, 17).Value = 25
$ws.Cells.Item(212, 18).Value = "Hortaliza"
